# Update column G ("K") values per regenerated save_data
# (commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    16 = 0
    19 = 1
    20 = 1
    21 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
